$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 69
$ws.Range("H69").Value = 4177
$ws.Range("I69").Value = 3980
$ws.Range("J69").Value = 4226.25
$ws.Range("K69").Value = 11940
$ws.Range("L69").Value = 12678.75
$ws.Range("M69").Value = -11066
$ws.Range("N69").Value = -14426.75
# Row 72
$ws.Range("H72").Value = 4177
$ws.Range("I72").Value = 3980
$ws.Range("J72").Value = 4226.25
$ws.Range("K72").Value = 35820
$ws.Range("L72").Value = 38036.25
$ws.Range("M72").Value = -31452
$ws.Range("N72").Value = -46772.25
# Row 76
$ws.Range("H76").Value = 2936.258
$ws.Range("I76").Value = 2698.5908
$ws.Range("J76").Value = 3517.2222
$ws.Range("K76").Value = 2698.5908
$ws.Range("L76").Value = 3517.2222
$ws.Range("M76").Value = -2383.5908
$ws.Range("N76").Value = -4147.2222
# Row 79
$ws.Range("H79").Value = 2936.258
$ws.Range("I79").Value = 2698.5908
$ws.Range("J79").Value = 3517.2222
$ws.Range("K79").Value = 2698.5908
$ws.Range("L79").Value = 3517.2222
$ws.Range("M79").Value = -1606.5908
$ws.Range("N79").Value = -5701.2222
# Row 100
$ws.Range("H100").Value = 3014.1667
$ws.Range("I100").Value = 2342.5
$ws.Range("K100").Value = 2342.5
$ws.Range("M100").Value = -1801.5
# Row 141
$ws.Range("H141").Value = 10408.625
$ws.Range("I141").Value = 2481.4443
$ws.Range("J141").Value = 20600.715
$ws.Range("K141").Value = 7444.3329
$ws.Range("L141").Value = 61802.145
$ws.Range("M141").Value = -2264.3329
$ws.Range("N141").Value = -72162.145

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("H97").Value = 274.66666
$ws.Range("I97").Value = 246.15384
$ws.Range("J97").Value = 460
$ws.Range("K97").Value = 246.15384
$ws.Range("L97").Value = 460
$ws.Range("M97").Value = 249.84616
$ws.Range("N97").Value = -1452

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 24400
$ws.Range("I82").Value = 2000
$ws.Range("K82").Value = 2000
$ws.Range("M82").Value = -1617
# Row 85
$ws.Range("H85").Value = 24400
$ws.Range("I85").Value = 2000
$ws.Range("K85").Value = 2000
$ws.Range("M85").Value = -674
# Row 94
$ws.Range("H94").Value = 327.21213
$ws.Range("I94").Value = 331.1875
$ws.Range("K94").Value = 331.1875
$ws.Range("M94").Value = 119.8125
# Row 99
$ws.Range("H99").Value = 2354.08
$ws.Range("I99").Value = 1965.2632
$ws.Range("J99").Value = 3585.3333
$ws.Range("K99").Value = 1965.2632
$ws.Range("L99").Value = 3585.3333
$ws.Range("M99").Value = -467.2632000000001
$ws.Range("N99").Value = -6581.3333

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 580
$ws.Range("I113").Value = 660
$ws.Range("J113").Value = 544.44446
$ws.Range("K113").Value = 1980
$ws.Range("L113").Value = 1633.33338
$ws.Range("M113").Value = 190
$ws.Range("N113").Value = -5973.33338

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3051.8572
$ws.Range("I80").Value = 2882.2222
$ws.Range("J80").Value = 3357.2
$ws.Range("K80").Value = 2882.2222
$ws.Range("L80").Value = 3357.2
$ws.Range("M80").Value = -1884.2222
$ws.Range("N80").Value = -5353.2
# Row 83
$ws.Range("H83").Value = 3051.8572
$ws.Range("I83").Value = 2882.2222
$ws.Range("J83").Value = 3357.2
$ws.Range("K83").Value = 14411.111
$ws.Range("L83").Value = 16786
$ws.Range("M83").Value = -9419.111
$ws.Range("N83").Value = -26770

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 6266.6665
$ws.Range("I2").Value = 4500
$ws.Range("J2").Value = 9800
$ws.Range("K2").Value = 4500
$ws.Range("L2").Value = 9800
$ws.Range("M2").Value = -4388
$ws.Range("N2").Value = -10024
# Row 3
$ws.Range("H3").Value = 10312.8125
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 10580.968
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 10580.968
$ws.Range("M3").Value = -1888
$ws.Range("N3").Value = -10804.968
# Row 10
$ws.Range("H10").Value = 10000000
$ws.Range("I10").Value = 10000000
$ws.Range("K10").Value = 10000000
$ws.Range("M10").Value = -9999860
# Row 12
$ws.Range("H12").Value = 2501325
$ws.Range("I12").Value = 5001900
$ws.Range("J12").Value = 750
$ws.Range("K12").Value = 5001900
$ws.Range("L12").Value = 750
$ws.Range("M12").Value = -5001730
$ws.Range("N12").Value = -1090
# Row 14
$ws.Range("H14").Value = 70005
$ws.Range("J14").Value = 70005
$ws.Range("L14").Value = 70005
$ws.Range("N14").Value = -70349
# Row 15
$ws.Range("H15").Value = 10312.8125
$ws.Range("I15").Value = 2000
$ws.Range("J15").Value = 10580.968
$ws.Range("K15").Value = 2000
$ws.Range("L15").Value = 10580.968
$ws.Range("M15").Value = -1830
$ws.Range("N15").Value = -10920.968
# Row 20
$ws.Range("H20").Value = 21999.3
$ws.Range("I20").Value = 23332.555
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 23332.555
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = -23106.555
$ws.Range("N20").Value = -10452
# Row 22
$ws.Range("H22").Value = 896.2692
$ws.Range("I22").Value = 829.64703
$ws.Range("J22").Value = 1022.1111
$ws.Range("K22").Value = 829.64703
$ws.Range("L22").Value = 1022.1111
$ws.Range("M22").Value = -534.64703
$ws.Range("N22").Value = -1612.1111
# Row 26
$ws.Range("H26").Value = 2004.5
$ws.Range("I26").Value = 2004.5
$ws.Range("K26").Value = 2004.5
$ws.Range("M26").Value = -1709.5
# Row 27
$ws.Range("H27").Value = 896.2692
$ws.Range("I27").Value = 829.64703
$ws.Range("J27").Value = 1022.1111
$ws.Range("K27").Value = 829.64703
$ws.Range("L27").Value = 1022.1111
$ws.Range("M27").Value = -722.64703
$ws.Range("N27").Value = -1236.1111
# Row 82
$ws.Range("H82").Value = 1778.8422
$ws.Range("I82").Value = 1370.3
$ws.Range("J82").Value = 2232.7778
$ws.Range("K82").Value = 1370.3
$ws.Range("L82").Value = 2232.7778
$ws.Range("M82").Value = -1009.3
$ws.Range("N82").Value = -2954.7778
# Row 85
$ws.Range("H85").Value = 1778.8422
$ws.Range("I85").Value = 1370.3
$ws.Range("J85").Value = 2232.7778
$ws.Range("K85").Value = 1370.3
$ws.Range("L85").Value = 2232.7778
$ws.Range("M85").Value = -122.3
$ws.Range("N85").Value = -4728.7778

